$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new portfolio snapshot row for 2025-09-24 (row 40), matching
# the existing layout: Date (text), SUZLON.NS, TATAMOTORS.NS, ETERNAL.NS.
$row = 40

# Force the date column to be stored as text so Excel doesn't silently
# reinterpret "2025-09-24" as a date serial number (all the other rows in
# column A are plain text, e.g. row 2..39).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-09-24"

$ws.Cells.Item($row, 2).Value = 57.5
$ws.Cells.Item($row, 3).Value = 682.9500122070312
$ws.Cells.Item($row, 4).Value = 335.7999877929688
